$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster data (Oyuncu Adı / Pozisyon / Takım), reordered with one
# player swapped out (Jared McCain -> Alexandre Sarr)
$players = @(
    @("Jalen Brunson",      "PG",       "New York Knicks"),
    @("Devin Booker",       "PG,SG",    "Phoenix Suns"),
    @("Desmond Bane",       "SG,SF",    "Memphis Grizzlies"),
    @("P.J. Washington",    "PF",       "Dallas Mavericks"),
    @("Alexandre Sarr",     "PF,C",     "Washington Wizards"),
    @("Yves Missi",         "C",        "New Orleans Pelicans"),
    @("Dereck Lively II",   "C",        "Dallas Mavericks"),
    @("LeBron James",       "SF,PF",    "Los Angeles Lakers"),
    @("Shaedon Sharpe",     "SG,SF",    "Portland Trail Blazers"),
    @("Jeremy Sochan",      "SF,PF",    "San Antonio Spurs"),
    @("Trae Young",         "PG",       "Atlanta Hawks"),
    @("Coby White",         "PG,SG",    "Chicago Bulls"),
    @("Walker Kessler",     "C",        "Utah Jazz"),
    @("Alperen Sengün",     "C",        "Houston Rockets"),
    @("Norman Powell",      "SG,SF",    "LA Clippers"),
    @("Immanuel Quickley",  "PG,SG",    "Toronto Raptors"),
    @("Kawhi Leonard",      "SG,SF,PF", "LA Clippers"),
    @("Devin Vassell",      "SG,SF",    "San Antonio Spurs")
)

for ($i = 0; $i -lt $players.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $players[$i][0]
    $ws.Cells.Item($row, 2).Value = $players[$i][1]
    $ws.Cells.Item($row, 3).Value = $players[$i][2]
}
